# "added first set of RF scores"
# Restructures the "Combined Results" and "Neural Network" sheets: re-orders /
# bolds the header row (with a bottom border), reshuffles the data rows, and
# fills in the first batch of Random Forest F1 scores on "Combined Results".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Combined Results"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Combined Results")

# Wipe the old data so we can lay the new layout down cleanly.
$ws1.Range("A1:G14").ClearContents()

# Header row (bold, thin bottom border).
$ws1.Range("A1").Value = "Dataset"
$ws1.Range("B1").Value = "Scaler"
$ws1.Range("C1").Value = "Features"
$ws1.Range("D1").Value = "Neural Network"
$ws1.Range("E1").Value = "Random Forest"
$ws1.Range("F1").Value = "Logistic Regression"
$ws1.Range("G1").Value = "Decision Tree"
$ws1.Range("A1:G1").Font.Bold = $true
$ws1.Range("A1:G1").Borders.Item(9).LineStyle = 1

# Data rows.
$ws1.Range("A2").Value = "Actual"
$ws1.Range("B2").Value = "Min-Max"
$ws1.Range("C2").Value = "All"

$ws1.Range("B3").Value = "Min-Max"
$ws1.Range("C3").Value = "Med Only"

$ws1.Range("B4").Value = "Standard"
$ws1.Range("C4").Value = "All"
$ws1.Range("E4").Value = "F1 = 0.80 / 0.82"

$ws1.Range("B5").Value = "Standard"
$ws1.Range("C5").Value = "Med Only"
$ws1.Range("E5").Value = "F1 = 0.82 /0.82"

$ws1.Range("B6").Value = "Unscaled"
$ws1.Range("C6").Value = "All"
$ws1.Range("E6").Value = "F1 = 0.80 / 0.82"

$ws1.Range("B7").Value = "Unscaled"
$ws1.Range("C7").Value = "Med Only"
$ws1.Range("E7").Value = "F1 = 0.83 / 0.83"

$ws1.Range("A8").Value = "RandomOverSampled"
$ws1.Range("B8").Value = "Min-Max"
$ws1.Range("C8").Value = "All"

$ws1.Range("B9").Value = "Min-Max"
$ws1.Range("C9").Value = "Med Only"

$ws1.Range("B10").Value = "Standard"
$ws1.Range("C10").Value = "All"

$ws1.Range("B11").Value = "Standard"
$ws1.Range("C11").Value = "Med Only"

$ws1.Range("A12").Value = "SMOTEENN"
$ws1.Range("B12").Value = "Min-Max"
$ws1.Range("C12").Value = "All"

$ws1.Range("B13").Value = "Min-Max"
$ws1.Range("C13").Value = "Med Only"

$ws1.Range("B14").Value = "Standard"
$ws1.Range("C14").Value = "All"

$ws1.Range("B15").Value = "Standard"
$ws1.Range("C15").Value = "Med Only"

$ws1.Range("D13").Select()

# ---------------------------------------------------------------------------
# Sheet 5: "Neural Network"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Neural Network")

$ws5.Range("A1:D14").ClearContents()

# Header row (bold, thin bottom border) - columns E:G stay blank but keep the
# header formatting, matching the "Combined Results" sheet's new layout.
$ws5.Range("A1").Value = "Dataset"
$ws5.Range("B1").Value = "Scaler"
$ws5.Range("C1").Value = "Features"
$ws5.Range("D1").Value = "Neural Network"
$ws5.Range("A1:G1").Font.Bold = $true
$ws5.Range("A1:G1").Borders.Item(9).LineStyle = 1

$ws5.Range("A2").Value = "Actual"
$ws5.Range("B2").Value = "Min-Max"
$ws5.Range("C2").Value = "All"

$ws5.Range("B3").Value = "Min-Max"
$ws5.Range("C3").Value = "Med Only"

$ws5.Range("B4").Value = "Standard"
$ws5.Range("C4").Value = "All"

$ws5.Range("B5").Value = "Standard"
$ws5.Range("C5").Value = "Med Only"

$ws5.Range("B6").Value = "Unscaled"
$ws5.Range("C6").Value = "All"

$ws5.Range("B7").Value = "Unscaled"
$ws5.Range("C7").Value = "Med Only"

$ws5.Range("A8").Value = "RandomOverSampled"
$ws5.Range("B8").Value = "Min-Max"
$ws5.Range("C8").Value = "All"

$ws5.Range("B9").Value = "Min-Max"
$ws5.Range("C9").Value = "Med Only"

$ws5.Range("B10").Value = "Standard"
$ws5.Range("C10").Value = "All"

$ws5.Range("B11").Value = "Standard"
$ws5.Range("C11").Value = "Med Only"

$ws5.Range("A12").Value = "SMOTEENN"
$ws5.Range("B12").Value = "Min-Max"
$ws5.Range("C12").Value = "All"

$ws5.Range("B13").Value = "Min-Max"
$ws5.Range("C13").Value = "Med Only"

$ws5.Range("B14").Value = "Standard"
$ws5.Range("C14").Value = "All"

$ws5.Range("B15").Value = "Standard"
$ws5.Range("C15").Value = "Med Only"

$ws5.PageSetup.Orientation = 1

$ws5.Range("E7").Select()

# ---------------------------------------------------------------------------
# Sheet 4: "Logistic Regression" - just a cursor-position change.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Logistic Regression")
$ws4.Range("D37").Select()

# ---------------------------------------------------------------------------
# Make "Combined Results" the active sheet/tab (was "Neural Network").
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D13").Select()
